# Updated capital structure database
# Refresh the Norway / Insurance (General) company rows (2-4) with the
# latest pulled figures: growth rates, margins, returns, cash/debt
# figures and the derived ratios that depend on them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (industry aggregate) ---------------------------------------
$ws.Range("D2").Value  = 0.10235
$ws.Range("E2").Value  = 0.0586
$ws.Range("F2").Value  = -0.08410000000000001
$ws.Range("G2").Value  = 0.2537134576261739
$ws.Range("H2").Value  = 0.2537134576261739
$ws.Range("I2").Value  = 0.1884179137753254
$ws.Range("J2").Value  = 0.1523075230380273
$ws.Range("K2").Value  = 541.3000000000001
$ws.Range("L2").Value  = 0.1508513780899033
$ws.Range("M2").Value  = 649.4200000000001
$ws.Range("N2").Value  = 0.05534279262005199
$ws.Range("O2").Value  = 1.199741363384445
$ws.Range("P2").Value  = 648.1
$ws.Range("Q2").Value  = 0.05523030380501939
$ws.Range("R2").Value  = 1.197302789580639
$ws.Range("S2").Value  = 1.32000000000005
$ws.Range("T2").Value  = 0.002032582920144205
$ws.Range("U2").Value  = 307.59
$ws.Range("V2").Value  = 0.02621245046657293
$ws.Range("W2").Value  = 0.2283490371435035
$ws.Range("X2").Value  = 0.04363630283617281
$ws.Range("Y2").Value  = 0.1847127343073307
$ws.Range("Z2").Value  = 1.318840047044987
$ws.Range("AA2").Value = 0.1757044250720492
$ws.Range("AB2").Value = 0.04085646650120132
$ws.Range("AC2").Value = 0.1348479585708479
$ws.Range("AD2").Value = 410.3
$ws.Range("AF2").Value = 410.3
$ws.Range("AG2").Value = 102.71
$ws.Range("AH2").Value = 0.0337840063236941
$ws.Range("AI2").Value = 0.1277755286350472
$ws.Range("AJ2").Value = 0.008676875716490624
$ws.Range("AK2").Value = 0.03537442612562036
$ws.Range("AN2").Value = 0.5819858156028369
$ws.Range("AP2").Value = 0.1456879432624113

# --- Row 3 (Gjensidige Forsikring ASA, OB:GJF) -------------------------
$ws.Range("D3").Value  = 0.0267
$ws.Range("E3").Value  = 0.0513
$ws.Range("F3").Value  = -0.08410000000000001
$ws.Range("G3").Value  = 0.3015653775322283
$ws.Range("H3").Value  = 0.3015653775322283
$ws.Range("I3").Value  = 0.2021507498026835
$ws.Range("J3").Value  = 0.1579189686924494
$ws.Range("K3").Value  = 480.1
$ws.Range("L3").Value  = 0.1578860826098395
$ws.Range("M3").Value  = 649.4200000000001
$ws.Range("N3").Value  = 0.05816517541266984
$ws.Range("O3").Value  = 1.352676525723808
$ws.Range("P3").Value  = 648.1
$ws.Range("Q3").Value  = 0.05804694987057886
$ws.Range("R3").Value  = 1.349927098521141
$ws.Range("S3").Value  = 1.32000000000005
$ws.Range("T3").Value  = 0.002032582920144205
$ws.Range("U3").Value  = 300.6
$ws.Range("V3").Value  = 0.02692318026708225
$ws.Range("W3").Value  = 0.175706338749817
$ws.Range("X3").Value  = 0.04110110972333884
$ws.Range("Y3").Value  = 0.1346052290264781
$ws.Range("Z3").Value  = 1.278990536277603
$ws.Range("AA3").Value = 0.2019768664563617
$ws.Range("AB3").Value = 0.04046753711324791
$ws.Range("AC3").Value = 0.1615093293431138
$ws.Range("AD3").Value = 277.5
$ws.Range("AF3").Value = 277.5
$ws.Range("AG3").Value = -23.10000000000002
$ws.Range("AH3").Value = 0.02425148130669603
$ws.Range("AI3").Value = 0.09898694442462723
$ws.Range("AJ3").Value = -0.002073236402800217
$ws.Range("AK3").Value = -0.009229662777688996
$ws.Range("AN3").Value = 0.4324450677886862
$ws.Range("AP3").Value = -0.03599812996727446

# --- Row 4 (Protector Forsikring ASA) ----------------------------------
# Ticker in the display name was corrected from OB:PROTCT to OB:PROT.
$ws.Range("B4").Value  = "Protector Forsikring ASA (OB:PROT)"
$ws.Range("D4").Value  = 0.178
$ws.Range("E4").Value  = 0.0659
$ws.Range("G4").Value  = -0.01205479452054794
$ws.Range("H4").Value  = -0.01205479452054794
$ws.Range("I4").Value  = 0.1121461187214612
$ws.Range("J4").Value  = 0.0936986301369863
$ws.Range("K4").Value  = 61.2
$ws.Range("L4").Value  = 0.1117808219178082
$ws.Range("M4").Value  = -0
$ws.Range("N4").Value  = -0
$ws.Range("O4").Value  = -0
$ws.Range("R4").Value  = -0
$ws.Range("S4").Value  = 0
$ws.Range("T4").ClearContents()
$ws.Range("U4").Value  = 6.99
$ws.Range("V4").Value  = 0.01227608008429926
$ws.Range("W4").Value  = 0.2809917355371901
$ws.Range("X4").Value  = 0.04617149594900677
$ws.Range("Y4").Value  = 0.2348202395881833
$ws.Range("Z4").Value  = 1.594815030585493
$ws.Range("AA4").Value = 0.1494319836877367
$ws.Range("AB4").Value = 0.04124539588915473
$ws.Range("AC4").Value = 0.1081865877985819
$ws.Range("AD4").Value = 132.8
$ws.Range("AF4").Value = 132.8
$ws.Range("AG4").Value = 125.81
$ws.Range("AH4").Value = 0.1891199088578752
$ws.Range("AI4").Value = 0.325729703213147
$ws.Range("AJ4").Value = 0.1809669020871391
$ws.Range("AK4").Value = 0.3139677073195079
$ws.Range("AN4").Value = 2.097946287519747
$ws.Range("AP4").Value = 1.987519747235387
